$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Coin/Link/Price/Volume columns touched below so Excel
# does not auto-coerce numeric-looking strings (e.g. "1.00", "28.00") into numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "98.691.49"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "3.336.71"
$ws.Range("E3").Value = "  +6.07%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "257.64"
$ws.Range("E5").Value = "  +6.44%  "

$ws.Range("D6").Value = "624.36"
$ws.Range("E6").Value = "  +2.45%  "

$ws.Range("D7").Value = "1.42"
$ws.Range("E7").Value = "  +27.83%  "

$ws.Range("D8").Value = "0.389"
$ws.Range("E8").Value = "  +2.00%  "

$ws.Range("E9").Value = "  -0.08%  "

$ws.Range("D10").Value = "0.878"
$ws.Range("E10").Value = "  +11.39%  "

$ws.Range("D11").Value = "3.334.96"
$ws.Range("E11").Value = "  +6.03%  "

$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "37.40"
$ws.Range("E13").Value = "  +10.31%  "

$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "98.345.53"
$ws.Range("E14").Value = "  +0.77%  "

$ws.Range("D15").Value = "0.0000248"
$ws.Range("E15").Value = "  +3.62%  "

$ws.Range("D16").Value = "3.944.92"
$ws.Range("E16").Value = "  +5.76%  "

$ws.Range("D17").Value = "5.52"
$ws.Range("E17").Value = "  +1.85%  "

$ws.Range("D18").Value = "3.329.66"
$ws.Range("E18").Value = "  +6.01%  "

$ws.Range("D19").Value = "3.54"
$ws.Range("E19").Value = "  +3.26%  "

$ws.Range("D20").Value = "15.14"
$ws.Range("E20").Value = "  +4.44%  "

$ws.Range("D21").Value = "489.44"
$ws.Range("E21").Value = "  -6.39%  "

$ws.Range("D22").Value = "6.10"
$ws.Range("E22").Value = "  +7.96%  "

$ws.Range("D23").Value = "0.0000210"
$ws.Range("E23").Value = "  +10.19%  "

$ws.Range("D24").Value = "9.35"
$ws.Range("E24").Value = "  +6.87%  "

$ws.Range("D25").Value = "5.61"
$ws.Range("E25").Value = "  +3.22%  "

$ws.Range("D26").Value = "88.61"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("D27").Value = "11.82"
$ws.Range("E27").Value = "  +2.16%  "

$ws.Range("D28").Value = "3.498.47"
$ws.Range("E28").Value = "  +5.55%  "

$ws.Range("D29").Value = "0.291"
$ws.Range("E29").Value = "  +22.74%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  +12.12%  "

$ws.Range("D32").Value = "0.138"
$ws.Range("E32").Value = "  +12.83%  "

$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("D34").Value = "9.65"
$ws.Range("E34").Value = "  +8.28%  "

$ws.Range("D35").Value = "28.00"
$ws.Range("E35").Value = "  +4.98%  "

$ws.Range("E36").Value = "  -0.87%  "

$ws.Range("D37").Value = "7.25"
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("D38").Value = "1.95"
$ws.Range("E38").Value = "  +4.04%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "497.48"
$ws.Range("E39").Value = "  +6.86%  "

$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "24.87"
$ws.Range("E40").Value = "  +2.11%  "

$ws.Range("D41").Value = "0.458"
$ws.Range("E41").Value = "  +5.54%  "

$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "1.26"
$ws.Range("E42").Value = "  +3.84%  "

$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").Value = "3.67"
$ws.Range("E43").Value = "  +4.74%  "

$ws.Range("D44").Value = "3.29"
$ws.Range("E44").Value = "  +6.10%  "

$ws.Range("D46").Value = "0.782"
$ws.Range("E46").Value = "  +12.62%  "

$ws.Range("D47").Value = "159.28"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  +1.56%  "

$ws.Range("D49").Value = "0.846"
$ws.Range("E49").Value = "  +9.68%  "

$ws.Range("D50").Value = "4.61"
$ws.Range("E50").Value = "  +3.22%  "

$ws.Range("D51").Value = "45.64"
$ws.Range("E51").Value = "  +3.65%  "
